{"js": "const replacements = [\n  ['79\u00f72=39, 1', '33\u00f74=8, 1'],\n  ['79\u00f78=9, 7', '18\u00f76=3, 0'],\n  ['19\u00f76=3, 1', '60\u00f72=30, 0'],\n  ['40\u00f78=5, 0', '73\u00f75=14, 3'],\n  ['35\u00f79=3, 8', '87\u00f72=43, 1'],\n  ['85\u00f72=42, 1', '78\u00f73=26, 0'],\n  ['71\u00f75=14, 1', '88\u00f77=12, 4'],\n  ['50\u00f79=5, 5', '62\u00f73=20, 2'],\n  ['13\u00f74=3, 1', '57\u00f74=14, 1'],\n  ['98\u00f72=49, 0', '55\u00f74=13, 3'],\n  ['60\u00f75=12, 0', '27\u00f76=4, 3'],\n  ['26\u00f75=5, 1', '92\u00f78=11, 4'],\n  ['38\u00f74=9, 2', '71\u00f78=8, 7'],\n  ['43\u00f74=10, 3', '57\u00f73=19, 0'],\n  ['15\u00f73=5, 0', '65\u00f76=10, 5'],\n  ['39\u00f73=13, 0', '35\u00f78=4, 3'],\n  ['63\u00f79=7, 0', '47\u00f72=23, 1'],\n  ['21\u00f78=2, 5', '18\u00f75=3, 3'],\n  ['33\u00f79=3, 6', '39\u00f78=4, 7'],\n  ['33\u00f75=6, 3', '35\u00f74=8, 3'],\n  ['34\u00f79=3, 7', '10\u00f77=1, 3'],\n  ['81\u00f79=9, 0', '69\u00f74=17, 1'],\n  ['76\u00f75=15, 1', '53\u00f74=13, 1'],\n  ['80\u00f77=11, 3', '40\u00f75=8, 0'],\n  ['94\u00f78=11, 6', '86\u00f73=28, 2'],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, 'Replace');\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"79\u00f72=39, 1\", \"33\u00f74=8, 1\"),\n    @(\"79\u00f78=9, 7\", \"18\u00f76=3, 0\"),\n    @(\"19\u00f76=3, 1\", \"60\u00f72=30, 0\"),\n    @(\"40\u00f78=5, 0\", \"73\u00f75=14, 3\"),\n    @(\"35\u00f79=3, 8\", \"87\u00f72=43, 1\"),\n    @(\"85\u00f72=42, 1\", \"78\u00f73=26, 0\"),\n    @(\"71\u00f75=14, 1\", \"88\u00f77=12, 4\"),\n    @(\"50\u00f79=5, 5\", \"62\u00f73=20, 2\"),\n    @(\"13\u00f74=3, 1\", \"57\u00f74=14, 1\"),\n    @(\"98\u00f72=49, 0\", \"55\u00f74=13, 3\"),\n    @(\"60\u00f75=12, 0\", \"27\u00f76=4, 3\"),\n    @(\"26\u00f75=5, 1\", \"92\u00f78=11, 4\"),\n    @(\"38\u00f74=9, 2\", \"71\u00f78=8, 7\"),\n    @(\"43\u00f74=10, 3\", \"57\u00f73=19, 0\"),\n    @(\"15\u00f73=5, 0\", \"65\u00f76=10, 5\"),\n    @(\"39\u00f73=13, 0\", \"35\u00f78=4, 3\"),\n    @(\"63\u00f79=7, 0\", \"47\u00f72=23, 1\"),\n    @(\"21\u00f78=2, 5\", \"18\u00f75=3, 3\"),\n    @(\"33\u00f79=3, 6\", \"39\u00f78=4, 7\"),\n    @(\"33\u00f75=6, 3\", \"35\u00f74=8, 3\"),\n    @(\"34\u00f79=3, 7\", \"10\u00f77=1, 3\"),\n    @(\"81\u00f79=9, 0\", \"69\u00f74=17, 1\"),\n    @(\"76\u00f75=15, 1\", \"53\u00f74=13, 1\"),\n    @(\"80\u00f77=11, 3\", \"40\u00f75=8, 0\"),\n    @(\"94\u00f78=11, 6\", \"86\u00f73=28, 2\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}"}
